$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1651.4193
$ws.Range("I15").Value = 1651.4193
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 4954.257900000001
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -4785.257900000001
# Row 17
$ws.Range("H17").Value = 5884927.5
$ws.Range("I17").Value = 1899
$ws.Range("J17").Value = 6252617
$ws.Range("K17").Value = 5697
$ws.Range("L17").Value = 18757851
$ws.Range("M17").Value = -5529
$ws.Range("N17").Value = -18758187
# Row 98
$ws.Range("H98").Value = 950.3333
$ws.Range("I98").Value = 913
$ws.Range("J98").Value = 1025
$ws.Range("K98").Value = 913
$ws.Range("L98").Value = 1025
$ws.Range("M98").Value = 585
$ws.Range("N98").Value = -4021
# Row 100
$ws.Range("H100").Value = 2034.1111
$ws.Range("I100").Value = 1546.7273
$ws.Range("J100").Value = 2800
$ws.Range("K100").Value = 1546.7273
$ws.Range("L100").Value = 2800
$ws.Range("M100").Value = -1005.7273
$ws.Range("N100").Value = -3882
# Row 112
$ws.Range("H112").Value = 1057.4073
$ws.Range("I112").Value = 800
$ws.Range("J112").Value = 1067.3077
$ws.Range("K112").Value = 2400
$ws.Range("L112").Value = 3201.9231
$ws.Range("M112").Value = -1292
$ws.Range("N112").Value = -5417.9231
# Row 116
$ws.Range("H116").Value = 4172.273
$ws.Range("I116").Value = 3222.25
$ws.Range("J116").Value = 4715.143
$ws.Range("K116").Value = 3222.25
$ws.Range("L116").Value = 4715.143
$ws.Range("M116").Value = 219.75
# Row 122
$ws.Range("H122").Value = 950.3333
$ws.Range("I122").Value = 913
$ws.Range("J122").Value = 1025
$ws.Range("K122").Value = 2739
$ws.Range("L122").Value = 3075
$ws.Range("M122").Value = -289
$ws.Range("N122").Value = -7975
# Row 129
$ws.Range("H129").Value = 179563.64
$ws.Range("I129").Value = 633.3333
$ws.Range("J129").Value = 189691.78
$ws.Range("K129").Value = 1899.9999
$ws.Range("L129").Value = 569075.34
$ws.Range("M129").Value = 3100.0001
$ws.Range("N129").Value = -579075.34

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 110
$ws.Range("H110").Value = 1027.75
$ws.Range("I110").Value = 1027.75
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1027.75
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1017.25
# Row 122
$ws.Range("H122").Value = 3402.875
$ws.Range("I122").Value = 2460.5715
$ws.Range("J122").Value = 9999
$ws.Range("K122").Value = 7381.7145
$ws.Range("L122").Value = 29997
$ws.Range("M122").Value = -4931.7145
$ws.Range("N122").Value = -34897

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2589.9546
$ws.Range("I20").Value = 2892.1333
$ws.Range("J20").Value = 1942.4286
$ws.Range("K20").Value = 2892.1333
$ws.Range("L20").Value = 1942.4286
$ws.Range("M20").Value = -2645.1333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 86
$ws.Range("H86").Value = 15434.214
$ws.Range("I86").Value = 7101.3
$ws.Range("J86").Value = 36266.5
$ws.Range("K86").Value = 7101.3
$ws.Range("L86").Value = 36266.5
$ws.Range("M86").Value = -5978.3
$ws.Range("N86").Value = -38512.5
# Row 89
$ws.Range("H89").Value = 15434.214
$ws.Range("I89").Value = 7101.3
$ws.Range("J89").Value = 36266.5
$ws.Range("K89").Value = 35506.5
$ws.Range("L89").Value = 181332.5
$ws.Range("M89").Value = -29890.5
$ws.Range("N89").Value = -192564.5
# Row 105
$ws.Range("H105").Value = 1074.6666
$ws.Range("I105").Value = 923.1667
$ws.Range("J105").Value = 1983.6666
$ws.Range("K105").Value = 923.1667
$ws.Range("L105").Value = 1983.6666
$ws.Range("M105").Value = 823.8333
$ws.Range("N105").Value = -5477.6666

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 14.8
$ws.Range("I2").Value = 15.692307
$ws.Range("J2").Value = 9
$ws.Range("K2").Value = 94.153842
$ws.Range("L2").Value = 54
$ws.Range("M2").Value = 18.846158
# Row 7
$ws.Range("H7").Value = 100
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 300
$ws.Range("N7").Value = -524
# Row 17
$ws.Range("H17").Value = 1088
$ws.Range("I17").Value = 326.66666
$ws.Range("J17").Value = 1414.2858
$ws.Range("K17").Value = 979.9999799999999
$ws.Range("L17").Value = 4242.857400000001
$ws.Range("M17").Value = -810.9999799999999
$ws.Range("N17").Value = -4580.857400000001
# Row 23
$ws.Range("H23").Value = 390.2
$ws.Range("I23").Value = 11
$ws.Range("J23").Value = 485
$ws.Range("K23").Value = 33
$ws.Range("L23").Value = 1455
$ws.Range("M23").Value = 202
$ws.Range("N23").Value = -1925
# Row 34
$ws.Range("H34").Value = 659.36365
$ws.Range("I34").Value = 112.25
$ws.Range("J34").Value = 972
$ws.Range("K34").Value = 336.75
$ws.Range("L34").Value = 2916
$ws.Range("M34").Value = -252.75
$ws.Range("N34").Value = -3084
# Row 55
$ws.Range("H55").Value = 2864.2856
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 2864.2856
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 8592.856800000001
$ws.Range("N55").Value = -8946.856800000001
# Row 103
$ws.Range("H103").Value = 1313.48
$ws.Range("I103").Value = 431.16666
$ws.Range("J103").Value = 2127.923
$ws.Range("K103").Value = 1293.49998
$ws.Range("L103").Value = 6383.768999999999
$ws.Range("M103").Value = -414.4999800000001
$ws.Range("N103").Value = -8141.768999999999
# Row 113
$ws.Range("H113").Value = 832.8
$ws.Range("I113").Value = 663.25
$ws.Range("J113").Value = 1026.5714
$ws.Range("K113").Value = 1989.75
$ws.Range("L113").Value = 3079.7142
$ws.Range("M113").Value = 180.25
$ws.Range("N113").Value = -7419.7142
# Row 131
$ws.Range("H131").Value = 722.7041
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 722.7041
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2168.1123
$ws.Range("N131").Value = -12248.1123
# Row 132
$ws.Range("H132").Value = 768.3333
$ws.Range("I132").Value = 700
$ws.Range("J132").Value = 802.5
$ws.Range("K132").Value = 6300
$ws.Range("L132").Value = 7222.5
$ws.Range("M132").Value = -3770
$ws.Range("N132").Value = -12282.5
# Row 136
$ws.Range("H136").Value = 3278.0715
# Row 138
$ws.Range("H138").Value = 2569.4707
$ws.Range("I138").Value = 2018.1818
$ws.Range("J138").Value = 3580.1667
$ws.Range("K138").Value = 6054.5454
$ws.Range("L138").Value = 10740.5001
$ws.Range("M138").Value = -914.5454
$ws.Range("N138").Value = -21020.5001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 7761.9546
$ws.Range("I70").Value = 3358.923
$ws.Range("J70").Value = 14121.889
$ws.Range("K70").Value = 3358.923
$ws.Range("L70").Value = 14121.889
$ws.Range("M70").Value = -3088.923
$ws.Range("N70").Value = -14661.889
# Row 73
$ws.Range("H73").Value = 7761.9546
$ws.Range("I73").Value = 3358.923
$ws.Range("J73").Value = 14121.889
$ws.Range("K73").Value = 3358.923
$ws.Range("L73").Value = 14121.889
$ws.Range("M73").Value = -2422.923
$ws.Range("N73").Value = -15993.889
# Row 92
$ws.Range("H92").Value = 14625
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 14625
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 14625
$ws.Range("N92").Value = -18369
$ws.Range("M92").ClearContents()
# Row 107
$ws.Range("H107").Value = 262
$ws.Range("I107").Value = 286.2
$ws.Range("J107").Value = 201.5
$ws.Range("K107").Value = 286.2
$ws.Range("L107").Value = 201.5
$ws.Range("M107").Value = 1633.8
$ws.Range("N107").Value = -4041.5
# Row 113
$ws.Range("H113").Value = 4842.7144
$ws.Range("I113").Value = 6626.8823
$ws.Range("J113").Value = 2085.3635
$ws.Range("K113").Value = 6626.8823
$ws.Range("L113").Value = 2085.3635
$ws.Range("M113").Value = -4456.8823
$ws.Range("N113").Value = -6425.363499999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 625.2
$ws.Range("I16").Value = 625.2
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 625.2
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -455.2
$ws.Range("N16").ClearContents()
# Row 34
$ws.Range("H34").Value = 80024
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 80024
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 80024
$ws.Range("N34").Value = -80368
# Row 35
$ws.Range("H35").Value = 168301.5
$ws.Range("I35").Value = 168301.5
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 168301.5
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -167965.5
# Row 68
$ws.Range("H68").Value = 2471.2856
$ws.Range("I68").Value = 1902
$ws.Range("J68").Value = 2566.1667
$ws.Range("K68").Value = 1902
$ws.Range("L68").Value = 2566.1667
$ws.Range("M68").Value = -1153
$ws.Range("N68").Value = -4064.1667
# Row 71
$ws.Range("H71").Value = 2471.2856
$ws.Range("I71").Value = 1902
$ws.Range("J71").Value = 2566.1667
$ws.Range("K71").Value = 9510
$ws.Range("L71").Value = 12830.8335
$ws.Range("M71").Value = -5766
$ws.Range("N71").Value = -20318.8335
# Row 93
$ws.Range("H93").Value = 1290.5714
$ws.Range("I93").Value = 1339
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 1339
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = -91
# Row 122
$ws.Range("H122").Value = 1786195.6
$ws.Range("I122").Value = 2453856.5
$ws.Range("J122").Value = 5766.6665
$ws.Range("K122").Value = 7361569.5
$ws.Range("L122").Value = 17299.9995
$ws.Range("M122").Value = -7359119.5
# Row 132
$ws.Range("H132").Value = 636319.6
$ws.Range("I132").Value = 1005316.25
$ws.Range("J132").Value = 3754
$ws.Range("K132").Value = 3015948.75
$ws.Range("L132").Value = 11262
$ws.Range("M132").Value = -3013418.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4687.875
$ws.Range("I62").Value = 3500
$ws.Range("J62").Value = 5083.8335
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 5083.8335
$ws.Range("M62").Value = -2876
# Row 65
$ws.Range("H65").Value = 4687.875
$ws.Range("I65").Value = 3500
$ws.Range("J65").Value = 5083.8335
$ws.Range("K65").Value = 17500
$ws.Range("L65").Value = 25419.1675
$ws.Range("M65").Value = -14380
# Row 96
$ws.Range("H96").Value = 1682.8572
$ws.Range("I96").Value = 1850
$ws.Range("J96").Value = 1616
$ws.Range("K96").Value = 1850
$ws.Range("L96").Value = 1616
$ws.Range("M96").Value = -477
$ws.Range("N96").Value = -4362
# Row 100
$ws.Range("H100").Value = 510
$ws.Range("I100").Value = 512
$ws.Range("J100").Value = 500
$ws.Range("K100").Value = 1024
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -483
$ws.Range("N100").Value = -2082
# Row 107
$ws.Range("H107").Value = 83334120
$ws.Range("I107").Value = 125000290
$ws.Range("J107").Value = 1775
$ws.Range("K107").Value = 375000870
$ws.Range("L107").Value = 5325
$ws.Range("M107").Value = -374998950
$ws.Range("N107").Value = -9165
# Row 132
$ws.Range("H132").Value = 1142.4667
$ws.Range("I132").Value = 824.9666999999999
$ws.Range("J132").Value = 1777.4667
$ws.Range("K132").Value = 2474.9001
$ws.Range("L132").Value = 5332.4001
$ws.Range("M132").Value = 55.09990000000016
$ws.Range("N132").Value = -10392.4001
